# GuildName.xlsx: unify the conception of DataNode, DataTable, Entity.
# The "Property" worksheet is renamed to "DataNode" to match the new naming
# convention, and the active selection is left where the editor last left it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet from "Property" to "DataNode".
$ws.Name = "DataNode"

# Leave the active selection on D39, matching the saved workbook view state.
[void]$ws.Range("D39").Select()
